$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("register")

# Add the new "promo code" related translation rows at the bottom of the
# register sheet (tag / english / dutch columns). Cells are populated in
# the exact order they were authored so that the shared-string table ends
# up in the same sequence as the source edit (the Dutch label for
# label_promo_code was filled in later, after the rows below it).

$ws.Cells.Item(62, 1).Value = "placeholder_promo_code"
$ws.Cells.Item(62, 2).Value = "promocode (optional)"
$ws.Cells.Item(62, 3).Value = "promocode (optioneel)"

$ws.Cells.Item(63, 1).Value = "label_promo_code"
$ws.Cells.Item(63, 2).Value = "Has this study been recommended to you by a friend? Fill in their promocode here (optional)"

$ws.Cells.Item(64, 1).Value = "invalid_feedback_promo_code"
$ws.Cells.Item(64, 2).Value = "This code does not exist"
$ws.Cells.Item(64, 3).Value = "Deze code bestaat niet"

$ws.Cells.Item(65, 1).Value = "valid_feedback_promo_code"
$ws.Cells.Item(65, 3).Value = "Deze code wordt gebruikt om degene die u heeft doorverwezen te belonen. Door de beloning zou hij/zij kunnen zien of u wel of niet de volledige periode deelneemt aan het onderzoek. Verder zullen er geen data van u worden gedeeld met degene die u heeft verwezen."
$ws.Cells.Item(65, 2).Value = "The promotion code is used for rewarding your acquaintance and thus he/she may know if you successfully participated or not. No other user data is shared with the acquaintance who gave you this code."

# The Dutch translation for "label_promo_code" (row trait) was added last.
$ws.Cells.Item(63, 3).Value = "Bent u doorverwezen naar dit onderzoek door een vriend? Vul hier hun promocode in (optioneel)"

# Reflect the scrolled/selected view state recorded after the edit.
$ws.Range("C63").Select()
$excel.ActiveWindow.ScrollRow = 53

$wb.Save()
